$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2022" column (S) mirroring the formatting already used for
#     the neighboring year columns in each row ---

# Header row (row 4): same look as the rest of the year headers (D4:Q4)
$ws.Range("D4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats - normalize R4 to the common header style
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Row 5 data
$ws.Range("E5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 49.7

# Row 6 data
$ws.Range("E6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 34.9

# Row 7 data (bottom border row)
$ws.Range("D7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 21

# --- Update the remembered selection to match the saved view ---
$ws.Range("R12").Select()
